# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# PIR sheet: append rows 255-265
# ---------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")

$pirRows = @(
    @("2026-02-01","14:14:57","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","14:14:59","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","14:15:04","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","14:15:09","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","14:15:14","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","14:15:19","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","14:15:24","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","14:15:29","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","14:15:34","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","14:15:38","14:00","Bathroom","Motion Detected","Active"),
    @("2026-02-01","14:15:42","14:00","Bathroom","Motion Detected","Active")
)

$startRow = 255
for ($i = 0; $i -lt $pirRows.Count; $i++) {
    $r = $startRow + $i
    $row = $pirRows[$i]

    $wsPIR.Cells.Item($r, 1).NumberFormat = "@"
    $wsPIR.Cells.Item($r, 1).Value = $row[0]

    $wsPIR.Cells.Item($r, 2).Value = $row[1]
    $wsPIR.Cells.Item($r, 3).Value = $row[2]
    $wsPIR.Cells.Item($r, 4).Value = $row[3]
    $wsPIR.Cells.Item($r, 5).Value = $row[4]
    $wsPIR.Cells.Item($r, 6).Value = $row[5]
}

# ---------------------------------------------------------------
# Humidity sheet: append rows 171-178
# ---------------------------------------------------------------
$wsHum = $wb.Worksheets.Item("Humidity")

$humRows = @(
    @("2026-02-01","14:14:56","14:00","Bathroom","77.9%","Active"),
    @("2026-02-01","14:14:58","14:00","Bathroom","76.9%","Active"),
    @("2026-02-01","14:15:07","14:00","Bathroom","76.9%","Active"),
    @("2026-02-01","14:15:12","14:00","Bathroom","77.9%","Active"),
    @("2026-02-01","14:15:22","14:00","Bathroom","77.0%","Active"),
    @("2026-02-01","14:15:32","14:00","Bathroom","77.0%","Active"),
    @("2026-02-01","14:15:42","14:00","Bathroom","77.0%","Active"),
    @("2026-02-01","14:15:52","14:00","Bathroom","78.7%","Active")
)

$startRow = 171
for ($i = 0; $i -lt $humRows.Count; $i++) {
    $r = $startRow + $i
    $row = $humRows[$i]

    $wsHum.Cells.Item($r, 1).NumberFormat = "@"
    $wsHum.Cells.Item($r, 1).Value = $row[0]

    $wsHum.Cells.Item($r, 2).Value = $row[1]
    $wsHum.Cells.Item($r, 3).Value = $row[2]
    $wsHum.Cells.Item($r, 4).Value = $row[3]

    $wsHum.Cells.Item($r, 5).NumberFormat = "@"
    $wsHum.Cells.Item($r, 5).Value = $row[4]

    $wsHum.Cells.Item($r, 6).Value = $row[5]
}

# ---------------------------------------------------------------
# Temperature sheet: append rows 92-99
# ---------------------------------------------------------------
$wsTemp = $wb.Worksheets.Item("Temperature")

$tempRows = @(
    @("2026-02-01","14:14:56","14:00","Bathroom","29.4C","Active"),
    @("2026-02-01","14:14:58","14:00","Bathroom","29.4C","Active"),
    @("2026-02-01","14:15:07","14:00","Bathroom","29.4C","Active"),
    @("2026-02-01","14:15:12","14:00","Bathroom","29.4C","Active"),
    @("2026-02-01","14:15:22","14:00","Bathroom","29.4C","Active"),
    @("2026-02-01","14:15:32","14:00","Bathroom","29.4C","Active"),
    @("2026-02-01","14:15:43","14:00","Bathroom","29.4C","Active"),
    @("2026-02-01","14:15:53","14:00","Bathroom","29.4C","Active")
)

$startRow = 92
for ($i = 0; $i -lt $tempRows.Count; $i++) {
    $r = $startRow + $i
    $row = $tempRows[$i]

    $wsTemp.Cells.Item($r, 1).NumberFormat = "@"
    $wsTemp.Cells.Item($r, 1).Value = $row[0]

    $wsTemp.Cells.Item($r, 2).Value = $row[1]
    $wsTemp.Cells.Item($r, 3).Value = $row[2]
    $wsTemp.Cells.Item($r, 4).Value = $row[3]
    $wsTemp.Cells.Item($r, 5).Value = $row[4]
    $wsTemp.Cells.Item($r, 6).Value = $row[5]
}
